$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 43, pushing the existing rows 43-55 down to 44-56.
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new weekly record.
$ws.Range("A43").Value = 4
$ws.Range("B43").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C43").Value = "Los Lagos"
$ws.Range("D43").Value = 44511
$ws.Range("E43").Value = 10
$ws.Range("F43").Value = 100112026
$ws.Range("G43").Value = "Haba"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 80
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = 10000
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Región Metropolitana"
$ws.Range("P43").Value = 400
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
